# "final checkin of day" - update the simulated sales figures on Blad1!B48:U48
# (row 48 is the raw random "times sold" input row; rows 49/50 and 53:58 are
# formulas that recompute automatically from it, and the chart on the
# "Chart2" chartsheet is fed from rows 53:58).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New values for B48:U48 (columns 2..21), replacing the old simulated data.
$newRow48 = @(1, 4, 6, 8, 4, 6, 2, 1, 3, 6, 7, 5, 9, 5, 4, 2, 3, 4, 6, 8)

for ($i = 0; $i -lt $newRow48.Length; $i++) {
    $ws.Cells.Item(48, 2 + $i).Value = $newRow48[$i]
}

# Match the author's final selection on this sheet.
$ws.Range("V48").Select() | Out-Null
